# "changing I1 to rerun spline"
# Sheet2 gets 24 new quarterly CET-1 ratio observations (2013Q3..2019Q2)
# inserted above the existing data, and the final row's quarter label is
# corrected from the erroneous "2024Q1" to "2023Q4".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert 24 blank rows above the existing data (rows 2:19 -> 26:43) and
# strip any inherited formatting so they stay unstyled like the rest of
# the data rows.
$ws2.Rows("2:25").Insert(1, 1)
$ws2.Range("A2:B25").ClearFormats()

$newDates = "2013Q3","2013Q4","2014Q1","2014Q2","2014Q3","2014Q4","2015Q1","2015Q2","2015Q3","2015Q4","2016Q1","2016Q2","2016Q3","2016Q4","2017Q1","2017Q2","2017Q3","2017Q4","2018Q1","2018Q2","2018Q3","2018Q4","2019Q1","2019Q2"
$newValues = 10.55,10.44,10.79,10.88,11.33,11.78,12.22,12.12,12.56,12.94,13.33,12.92,13.25,12.89,12.9,12.56,12.62,12.87,12.98,12.87,13.17,13.74,13.71,13.47

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws2.Range("A$row").Value = $newDates[$i]
    $ws2.Range("B$row").Value = $newValues[$i]
}

# The old last row (originally A19/B19, now shifted to row 43) had a stray
# "2024Q1" label -- correct it to "2023Q4" now that all the new quarter
# strings above have already been introduced, so the shared-string table
# ends up with the new quarters first and "2023Q4" last, matching the
# target layout once the now-unused "2024Q1" entry is compacted away.
$ws2.Range("A43").Value = "2023Q4"

# B1's header cell loses its grey shading (A1 keeps it).
$ws2.Range("B1").ClearFormats()

# --- sheet selections / active sheet -----------------------------------

$ws1.Activate()
$ws1.Range("C551").Select()

$ws2.Activate()
$ws2.Range("A44").Select()
